# Apply "Horarios actualizados Linea 141 - 133" update
# Updates header metadata (timestamp / row counts) and schedule rows
# across the three worksheets: LP1912, LP1912-215, 6203-6173.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 06:54:04"
$ws1.Range("A3").Value = "Total filas: 39"

$sheet1Rows = @(
    @("06:54:04","07:22","23_HERNANDEZ",28,"LP1912"),
    @("05:57:13","07:29","17X38_ROMERO",92,"LP1912"),
    @("05:57:13","07:35","10_OLMOS",98,"LP1912"),
    @("06:17:28","07:36","27_EL RETIRO",79,"LP1912"),
    @("05:57:13","07:37","27_EL RETIRO",100,"LP1912"),
    @("06:46:50","07:43","215A_EL PATO",57,"LP1912"),
    @("06:35:22","07:44","215A_EL PATO",69,"LP1912"),
    @("05:57:13","07:55","14_ABASTO",118,"LP1912"),
    @("06:17:28","08:00","17_ROMERO",103,"LP1912"),
    @("06:46:50","08:00","16_SANTA ANA",74,"LP1912"),
    @("06:17:28","08:01","16_SANTA ANA",104,"LP1912"),
    @("06:35:22","08:06","23_HERNANDEZ",91,"LP1912"),
    @("06:54:04","08:07","23_HERNANDEZ",73,"LP1912"),
    @("06:17:28","08:11","10_OLMOS",114,"LP1912"),
    @("06:17:28","08:13","15X38_ABASTO",116,"LP1912"),
    @("06:35:22","08:29","11_ETCHEVERRY",114,"LP1912"),
    @("06:35:22","08:29","15_ABASTO",114,"LP1912"),
    @("06:46:50","08:41","16_P MOR-SANTA ANA",115,"LP1912"),
    @("06:46:50","08:43","215C_EL PATO",117,"LP1912"),
    @("06:54:04","08:44","215C_EL PATO",110,"LP1912")
)

$startRow = 25
for ($i = 0; $i -lt $sheet1Rows.Count; $i++) {
    $r = $startRow + $i
    $row = $sheet1Rows[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
}

# ---------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 06:54:04"
$ws2.Range("A3").Value = "Total filas: 10"

$ws2.Cells.Item(15, 1).Value = "06:54:04"
$ws2.Cells.Item(15, 2).Value = "08:44"
$ws2.Cells.Item(15, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(15, 4).Value = 110
$ws2.Cells.Item(15, 5).Value = "LP1912"

# ---------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 06:54:04"
$ws3.Range("A3").Value = "Total filas: 5"

$ws3.Cells.Item(9, 1).Value = "06:54:04"
$ws3.Cells.Item(9, 2).Value = "08:36"
$ws3.Cells.Item(9, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(9, 4).Value = 102
$ws3.Cells.Item(9, 5).Value = "L6173"

$ws3.Cells.Item(10, 1).Value = "06:54:04"
$ws3.Cells.Item(10, 2).Value = "08:51"
$ws3.Cells.Item(10, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(10, 4).Value = 117
$ws3.Cells.Item(10, 5).Value = "L6203"

Write-Output "Horarios actualizados Linea 141 - 133"
